$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title of the first (Books) table ---
$ws.Range("A1").Value = "Table for Books"

# --- Header row (row 2) relabeling ---
$ws.Range("G2").Value = "Image source"
$ws.Range("K2").Value = "No. of veiws"
$ws.Range("L2").Value = "InStock"
$ws.Range("M2").Value = "NumSold"

# --- Data type row (row 3) updates ---
$ws.Range("K3").Value = "No. of veiws"
$ws.Range("L3").Value = "boolean"
$ws.Range("M3").Value = "int"

# --- Remove the now-unused Purchase History table (rows 10-12) ---
$ws.Range("A10:N12").Delete()

# --- Remove the now-unused trailing column N (previously NumSold) ---
$ws.Columns.Item(14).Delete()

# --- Column width / formatting tweaks ---
$ws.Columns.Item(1).ColumnWidth = 24.1
$ws.Columns.Item(11).ColumnWidth = 15.45

# --- Selection moves to H32 ---
[void]$ws.Range("H32").Select()
